# GanttChart.xlsx update:
# - Merge 3 pairs of tasks ("Create Basic X" + "Update ... to make/return ...")
#   into single combined task descriptions (rows 4-6 on the "Gantt Chart" sheet).
# - Remove the now-unused trailing rows (old rows 23-29), which drops the last
#   caching task row and the un-scheduled voice-recognition / machine-learning
#   rows that previously had no dates (they are now scheduled in rows 17-22).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gantt Chart")

# Delete the trailing rows 23:29 (shifts rows below up - none exist below row 29).
$ws.Rows("23:29").Delete()

# Re-write the task names (column A) for rows 4 through 22 with the
# consolidated / re-ordered descriptions.
$tasks = @(
    "Create Basic Client Program to make web requests",
    "Create Basic Edge Server Program to return web response",
    "Create Data Centre program to return web responses",
    "Update Client Program Design",
    "Update Edge Server Program to receive web requests",
    "Update Edge Server Program to process web requests",
    "Update Data Centre Program to receive web requests",
    "Update Data Centre program to process web requests",
    "Update Client Program to handle responses",
    "Update Client Program to make multiple web requests",
    "Research and design caching application ",
    "Create image for the cachine application",
    "Perform analysis and update caching application",
    "Research and design voice recognition application",
    "Create image for voice recognition application",
    "Perform analysis and update voice recognition application",
    "Research and design machine learning application",
    "Create image for machine learning application",
    "Perform analysis and update machine learning application"
)

for ($i = 0; $i -lt $tasks.Count; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 1).Value = $tasks[$i]
}
